# Logged Week 15 and simulated Week 16
$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# YDS sheet - append the per-play yardage logs for the newly logged /
# simulated weeks to the running totals.
# ----------------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")

$ydsWs.Range("B2").Value = $ydsWs.Range("B2").Value2 + " 4 15 3 15 4 3 6 2 4 1 2 7 0 9 4 6 8 4 3 1 8 2 4 8 22 1 1 3 0 2 2 3 2 6 15 2 2"
$ydsWs.Range("B3").Value = $ydsWs.Range("B3").Value2 + " 9 4 12 12 13 15 4 2 27 7 5 4 17 6 29 8 5 15 15 8 5 14"
$ydsWs.Range("C2").Value = $ydsWs.Range("C2").Value2 + " 5 11 9 7 0 5 0 -1 3 3 1 2 3 0 1 5 0 1 32"
$ydsWs.Range("C3").Value = $ydsWs.Range("C3").Value2 + " 12 14 5 9 20 11 12 -1 17 4 3 6 5 23 10 40 5 3 9 0 6 69 1 5 10 20 7 10 27 14 34"

# ----------------------------------------------------------------------
# OFF sheet - updated season totals
# ----------------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")

$offWs.Range("C2").Value = 162
$offWs.Range("D2").Value = 15
$offWs.Range("F2").Value = 42
$offWs.Range("G2").Value = 50
$offWs.Range("J2").Value = 18
$offWs.Range("L2").Value = 303
$offWs.Range("M2").Value = 195
$offWs.Range("O2").Value = 26
$offWs.Range("P2").Value = 15
$offWs.Range("Q2").Value = 528

$offWs.Range("C3").Value = 175
$offWs.Range("D3").Value = 5
$offWs.Range("E3").Value = 34
$offWs.Range("F3").Value = 98
$offWs.Range("G3").Value = 43
$offWs.Range("H3").Value = 27
$offWs.Range("I3").Value = 52
$offWs.Range("J3").Value = 55

# ----------------------------------------------------------------------
# DEF sheet - updated season totals
# ----------------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")

$defWs.Range("B2").Value = 3
$defWs.Range("C2").Value = 193
$defWs.Range("D2").Value = 13
$defWs.Range("F2").Value = 51
$defWs.Range("G2").Value = 65
$defWs.Range("J2").Value = 36
$defWs.Range("L2").Value = 293
$defWs.Range("M2").Value = 186
$defWs.Range("O2").Value = 17
$defWs.Range("Q2").Value = 528

$defWs.Range("C3").Value = 153
$defWs.Range("D3").Value = 4
$defWs.Range("E3").Value = 31
$defWs.Range("F3").Value = 83
$defWs.Range("G3").Value = 43
$defWs.Range("I3").Value = 46
$defWs.Range("J3").Value = 48
$defWs.Range("N3").Value = 16

# ----------------------------------------------------------------------
# ST sheet - updated season totals and appended per-week logs
# ----------------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")

$stWs.Range("B2").Value = 78
$stWs.Range("D2").Value = 42
$stWs.Range("F2").Value = 218
$stWs.Range("G2").Value = 205

$stWs.Range("B3").Value = 46

$stWs.Range("B4").Value = $stWs.Range("B4").Value2 + " 69"
$stWs.Range("D3").Value = $stWs.Range("D3").Value2 + " 41"
$stWs.Range("D4").Value = $stWs.Range("D4").Value2 + " 8"
$stWs.Range("B5").Value = $stWs.Range("B5").Value2 + " 19"
$stWs.Range("D5").Value = $stWs.Range("D5").Value2 + " 0"
$stWs.Range("B6").Value = $stWs.Range("B6").Value2 + " 75"

# ----------------------------------------------------------------------
# TURNS sheet - updated turnover totals
# ----------------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")

$turnsWs.Range("B2").Value = 7
$turnsWs.Range("C2").Value = 4
$turnsWs.Range("D2").Value = 3
$turnsWs.Range("E2").Value = 8

$turnsWs.Range("E3").Value = 7

# ----------------------------------------------------------------------
# PEN sheet - updated penalty totals
# ----------------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")

$penWs.Range("D4").Value = 7
